$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 93

$ws.Cells.Item($row, 1).Value = "Record"
$ws.Cells.Item($row, 2).Value = "RJ Record"
$ws.Cells.Item($row, 3).Value = "Agricultura"
$ws.Cells.Item($row, 4).Value = "2025-04-08T18:53"
$ws.Cells.Item($row, 5).Value = "Positivo"
$ws.Cells.Item($row, 6).Value = "Maior cidade do interior do Rio se torna principal produtor de soja do Estado. Entrevista com produtor rural José Geraldo Neto. Campos vai exportar 1.800 tonadas e Macaé 1.100 toneladas, totalizando 3 mil toneladas. Segundo dados da Embrapa, são cerca de 300 mil hectares com aptidão para a soja. Proximidade com o Porto do Açu. Essa é a sexta vez que o produto será enviado a Rússia pelo Porto do Açu. Entrevista com o secretário de Agricultura, Almy Junior, e com gerente de propriedade Manuel Peixoto.  "
